# Add three new "cht*_f_e1" helper-data sheets after cht11_f_e1, each
# holding the same small relative-depth/f_e1 lookup bounds (1,1) and
# (100,100) used to let downstream LOOKUP/INTERP formulas be cached
# instead of re-read from the big cht11_f_e1 table every time.

$wb = $excel.ActiveWorkbook

$sheetNames = @("cht12_f_e1", "cht13_f_e1", "cht14_f_e1")

foreach ($name in $sheetNames) {
    $afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $afterSheet)
    $ws.Name = $name

    $ws.Range("A1").Value = "e_ss"
    $ws.Range("B1").Value = "f_e1"
    $ws.Range("A2").Value = 1
    $ws.Range("B2").Value = 1
    $ws.Range("A3").Value = 100
    $ws.Range("B3").Value = 100
}

# Restore per-sheet selections like a normal editing session would leave
# behind, then land back on cht12_f_e1 as the active tab.
$wb.Worksheets.Item("cht13_f_e1").Range("B4").Select() | Out-Null
$wb.Worksheets.Item("cht14_f_e1").Range("B4").Select() | Out-Null

$ws12 = $wb.Worksheets.Item("cht12_f_e1")
$ws12.Activate()
$ws12.Rows("4:72").Select() | Out-Null
